# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B, C, D, E, F stays same, G = B+C+D+E)
$data = @{
    2 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    3 = @(0.0008583669626518464, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0.6975673568529668)
    4 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 8.660232485948974, 38.00504557982321)
    5 = @(0.3048080303191223, 0.3127903958511391, 3.900430680208489, 8.660232485948974, 13.17826159232772)
    6 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    7 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 17.45944343273191)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
